$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections (reverse-complement / rotation clade fixes) ---
# KJ641714 / BtRp-CV-14 was misclassified as genus "Cyclovirus" with clade
# "Cyclovirus-3"; it is actually "Circovirus" with no clade assigned.
$ws.Range("E46").Value = "Circovirus"
$ws.Range("F46").Value = ""

# Re-sort the data block (rows 2:49) by genus (E) then clade (F), ascending,
# matching the sheet's existing sortState. This rotates the corrected
# KJ641714 row up to the top of its genus group (blanks sort last) and
# re-groups the Cyclovirus/CRESS blocks accordingly.
$rng = $ws.Range("A2:H49")
$keyE = $ws.Range("E2:E49")
$keyF = $ws.Range("F2:F49")
$rng.Sort($keyE, 1, $keyF, $null, 1, $null, $null, 1)

# PCV-3 (NC_031753) clade correction: Mammal-2 -> Mammal-1
$ws.Range("F20").Value = "Mammal-1"
